$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are stored as plain text in this sheet (prices/volumes are
# formatted strings like "67.693.11" or "  -1.21%  "). For values that would
# otherwise be auto-recognized by Excel as a number, force the cell to text
# first so the literal string is preserved, then clear the temporary format
# so no stray NumberFormat/style is left behind.

$ws.Range("D2").Value = "67.693.11"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "3.773.83"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "595.04"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.17%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "167.37"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "3.770.82"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -0.56%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.31"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -2.21%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.447"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  -2.64%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "36.06"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "4.406.50"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "3.810.64"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "67.641.09"
$ws.Range("E17").Value = "  -1.33%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "18.28"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +1.62%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.97"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.68%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.01"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -6.88%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "456.49"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -1.72%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.694"
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.0000153"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("E25").Value = "  -1.24%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.90"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.28%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.13"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -0.21%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.22"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("E32").Value = "  -0.77%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "29.62"
$cell.ClearFormats()
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "9.14"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "3.726.67"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  -0.25%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "3.33"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("E39").Value = "  -1.06%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -0.41%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.75"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "46.50"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +5.66%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -0.02%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "48.24"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +3.62%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.299"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.13%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "149.01"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("E48").Value = "  -2.10%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "389.27"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.54%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.82"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -4.94%  "
$ws.Range("E51").Value = "  +0.02%  "
